# Swap the taxon/observation data between row 2 and row 3 in the
# "Artfynd" sheet. Columns A, B, D, E, F, G, H, Q, R hold the data that
# needs to be exchanged between the two rows; the remaining columns are
# identical between the rows (shared location/date/observer info) and
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")

    $val2 = $cell2.Value2
    $val3 = $cell3.Value2

    $cell2.Value2 = $val3
    $cell3.Value2 = $val2
}
